$wb = $excel.ActiveWorkbook

# Work on the "Products" worksheet (sheet2.xml)
$ws = $wb.Worksheets.Item("Products")

$ws.Range("A1").Value = "Iphone 12"
$ws.Range("A2").Value = "Mobile"
$ws.Range("A3").Value = 112000

$ws.Columns.Item(1).ColumnWidth = 17.65

$ws.Activate()
$ws.Range("A3").Select()
